# Correct output files with 2018 and trying to add lines to plots
#
# Insert a new row of data for LOC_2018 before the existing LOC_2019 row
# (which currently sits at row 6), shifting LOC_2019..LOC_2023 down by one
# row. The dimension grows from A1:J10 to A1:J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 6 (this pushes the old row 6 -> 7, etc.)
$ws.Rows.Item(6).Insert()

# Force the new row's cells to be stored as text, matching the rest of the
# sheet where every value (including numbers) is written as a text string.
$newRow = $ws.Range("A6:J6")
$newRow.NumberFormat = "@"

$ws.Range("A6").Value2 = "LOC_2018"
$ws.Range("B6").Value2 = "2018-05-08"
$ws.Range("C6").Value2 = "2018-05-21"
$ws.Range("D6").Value2 = "282.68"
$ws.Range("E6").Value2 = "449.07"
$ws.Range("F6").Value2 = "0.99792088406833"
$ws.Range("G6").Value2 = "1.823776414708e-17"
$ws.Range("H6").Value2 = "0.00015377009564513"
$ws.Range("I6").Value2 = "-234335.931116353"
$ws.Range("J6").Value2 = "functional_ice_off_to_no_ice"
